$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (target OOXML width is 32.140625 characters)
$ws.Columns.Item(1).ColumnWidth = 31.33

# Replace the scheduled-time numeric date value in A2 with a text
# representation of the date/time (keeps the cell's existing style,
# simply switches the stored value to a shared string).
$ws.Range("A2").Value = "2017-06-07T13:34:08.0039447-05:00"
